$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "30.040.35"
Set-TextValue "E2" "  +7.61%  "
Set-TextValue "D3" "1.874.78"
Set-TextValue "E3" "  +5.37%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "248.33"
Set-TextValue "E5" "  +2.06%  "
Set-TextValue "D6" "0.9998"
Set-TextValue "E6" "  -0.01%  "
Set-TextValue "D7" "0.4963"
Set-TextValue "E7" "  +1.53%  "
Set-TextValue "D8" "45.72"
Set-TextValue "E8" "  +9.05%  "
Set-TextValue "E9" "  +6.89%  "
Set-TextValue "D10" "0.06559"
Set-TextValue "E10" "  +5.10%  "
Set-TextValue "D11" "1.872.50"
Set-TextValue "E11" "  +5.27%  "
Set-TextValue "D12" "17.06"
Set-TextValue "E12" "  +4.84%  "
Set-TextValue "D13" "0.07183"
Set-TextValue "E13" "  +2.63%  "
Set-TextValue "D14" "0.6606"
Set-TextValue "E14" "  +6.87%  "
Set-TextValue "D15" "85.09"
Set-TextValue "E15" "  +6.86%  "
Set-TextValue "D16" "4.791"
Set-TextValue "E16" "  +3.80%  "
Set-TextValue "D17" "30.010.05"
Set-TextValue "E17" "  +7.62%  "
Set-TextValue "D19" "12.82"
Set-TextValue "E19" "  +8.55%  "
Set-TextValue "D20" "0.000007488"
Set-TextValue "E20" "  +4.06%  "
Set-TextValue "D21" "0.9999"
Set-TextValue "E21" "  -0.01%  "
Set-TextValue "D22" "2.113.76"
Set-TextValue "E22" "  +5.23%  "
Set-TextValue "D23" "4.728"
Set-TextValue "E23" "  +3.59%  "
Set-TextValue "D24" "9.022"
Set-TextValue "E24" "  +4.30%  "
Set-TextValue "E25" "  +5.64%  "
Set-TextValue "D26" "144.09"
Set-TextValue "E26" "  +1.83%  "
Set-TextValue "D27" "135.13"
Set-TextValue "E27" "  +24.18%  "
Set-TextValue "D28" "16.69"
Set-TextValue "E28" "  +7.28%  "
Set-TextValue "D29" "1.954"
Set-TextValue "E29" "  +5.25%  "
Set-TextValue "D30" "1.385"
Set-TextValue "E30" "  -0.55%  "
Set-TextValue "D31" "4.204"
Set-TextValue "E31" "  +1.37%  "
Set-TextValue "E32" "  +4.34%  "
Set-TextValue "D33" "3.878"
Set-TextValue "E33" "  +2.44%  "
Set-TextValue "D34" "0.05062"
Set-TextValue "E34" "  +6.51%  "
Set-TextValue "E35" "  +6.58%  "
Set-TextValue "D36" "0.6819"
Set-TextValue "E36" "  +6.04%  "
Set-TextValue "D37" "2.702"
Set-TextValue "E37" "  +3.86%  "
Set-TextValue "D38" "2.309"
Set-TextValue "E38" "  +13.36%  "
Set-TextValue "D39" "2.741"
Set-TextValue "E39" "  +5.85%  "
Set-TextValue "D40" "0.9621"
Set-TextValue "E40" "  +2.26%  "
Set-TextValue "E41" "  +5.80%  "
Set-TextValue "D42" "6.054"
Set-TextValue "E42" "  +2.71%  "
Set-TextValue "D43" "1.000"
Set-TextValue "E43" "  -0.01%  "
Set-TextValue "D44" "103.13"
Set-TextValue "E44" "  +3.12%  "
Set-TextValue "D45" "0.4188"
Set-TextValue "E45" "  +6.03%  "
Set-TextValue "D46" "7.446"
Set-TextValue "E46" "  +3.13%  "
Set-TextValue "D47" "0.1254"
Set-TextValue "D48" "0.05622"
Set-TextValue "E48" "  +3.88%  "
Set-TextValue "D49" "32.36"
Set-TextValue "E49" "  +6.43%  "
Set-TextValue "D50" "8.238"
Set-TextValue "E50" "  +3.13%  "
Set-TextValue "D51" "0.3720"
Set-TextValue "E51" "  +7.56%  "
